# Update on 20250925 part 5
# Add a new "其他备注" (Other Notes) worksheet at the end of the workbook
# documenting proxy-source <-> channel mappings for a few sports channels.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet (so it lands at the end,
# matching the workbook.xml order: 央视, CETV和CGTN, 上海, 卫视, 体育, 娱乐, 其他备注)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "其他备注"

# Colors (Excel BGR-packed long values)
$yellow = 65535      # RGB FFFF00
$blue   = 15773696   # RGB 00B0F0

# ---- Cells: address, text value (or $null for a blank-but-bordered cell),
#      and interior fill color (or $null for no fill) ---------------------
$cells = @(
  @{ Addr = "A1"; Value = $null;              Fill = $null  }
  @{ Addr = "B1"; Value = "当前源";            Fill = $null  }
  @{ Addr = "C1"; Value = "代理源失效,可平替观赛频道"; Fill = $null  }
  @{ Addr = "E1"; Value = "说明：";            Fill = $null  }
  @{ Addr = "F1"; Value = $null;              Fill = $null  }

  @{ Addr = "A2"; Value = "CCTV-5+";          Fill = $null  }
  @{ Addr = "B2"; Value = "咪咕代理源";        Fill = $yellow }
  @{ Addr = "C2"; Value = "咪视界";            Fill = $blue   }
  @{ Addr = "E2"; Value = $null;              Fill = $yellow }
  @{ Addr = "F2"; Value = "为不稳定源";        Fill = $null  }

  @{ Addr = "A3"; Value = "劲爆体育";          Fill = $null  }
  @{ Addr = "B3"; Value = "咪咕代理源";        Fill = $yellow }
  @{ Addr = "C3"; Value = "咪视界";            Fill = $blue   }
  @{ Addr = "E3"; Value = $null;              Fill = $blue   }
  @{ Addr = "F3"; Value = "为移动网络专用源";  Fill = $null  }

  @{ Addr = "A4"; Value = "睛彩青少";          Fill = $null  }
  @{ Addr = "B4"; Value = "肥羊代理源";        Fill = $yellow }
  @{ Addr = "C4"; Value = "咪视界";            Fill = $blue   }
)

foreach ($cell in $cells) {
  $rng = $ws.Range($cell.Addr)
  if ($cell.Value -ne $null) {
    $rng.Value2 = $cell.Value
  }
  $rng.Borders.LineStyle = 1
  $rng.Borders.Weight = 2
  if ($cell.Fill -ne $null) {
    $rng.Interior.Color = $cell.Fill
  }
}

# ---- Column widths (characters) - closest achievable values to the
#      authored 11 / 26.625 / 7.125 / 17.25 widths ------------------------
$ws.Columns.Item(2).ColumnWidth = 10.285714285714286
$ws.Columns.Item(3).ColumnWidth = 25.857142857142858
$ws.Columns.Item(5).ColumnWidth = 6.428571428571429
$ws.Columns.Item(6).ColumnWidth = 16.571428571428573
